# Updated cryptos list on Sat Feb  3 06:53:31 UTC 2024 with GitHub Actions
#
# Applies the per-row Price / Volume(1h) refresh to the crypto tracker sheet,
# plus the rank swap between BitcoinSV and MultiversX (rows 48/49).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell while forcing it to stay text (Excel's
# COM layer auto-coerces plain decimal-looking strings like "303.57" or
# "1.00" into numbers on a bare .Value assignment, which would silently
# drop the fixed 2-decimal display formatting the source data relies on).
# Flip the cell to text ("@") first, assign, then restore a styleless
# ("Normal") cell so no stray number-format sticks around on the cell.
function Set-TextValue {
    param($cell, [string]$value)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# row -> @{ B=...; C=...; D=...; E=... } (only the columns that changed)
$rowUpdates = [ordered]@{
    2  = @{ D = "43.114.12";  E = "  -0.29%  " }
    3  = @{ D = "2.320.54";   E = "  +0.41%  " }
    4  = @{ E = "  +0.05%  " }
    5  = @{ D = "303.57";     E = "  +0.34%  " }
    6  = @{ D = "99.74";      E = "  -0.30%  " }
    7  = @{ E = "  +0.64%  " }
    9  = @{ D = "0.520";      E = "  +2.19%  " }
    10 = @{ D = "36.09";      E = "  +4.53%  " }
    11 = @{ E = "  -0.93%  " }
    12 = @{ E = "  -1.56%  " }
    13 = @{ D = "17.79";      E = "  -1.53%  " }
    14 = @{ D = "6.93";       E = "  +1.23%  " }
    15 = @{ D = "2.681.58";   E = "  +0.52%  " }
    16 = @{ D = "2.312.48";   E = "  +1.43%  " }
    17 = @{ E = "  -2.84%  " }
    18 = @{ D = "43.034.49";  E = "  -0.28%  " }
    19 = @{ D = "13.25";      E = "  +4.38%  " }
    20 = @{ D = "6.19";       E = "  +1.12%  " }
    21 = @{ E = "  +0.17%  " }
    22 = @{ E = "  +0.57%  " }
    23 = @{ D = "239.95";     E = "  +0.94%  " }
    24 = @{ E = "  -2.04%  " }
    25 = @{ D = "2.47";       E = "  -0.09%  " }
    26 = @{ D = "1.00";       E = "  -0.02%  " }
    27 = @{ D = "25.44";      E = "  +2.45%  " }
    28 = @{ D = "169.15";     E = "  +0.30%  " }
    29 = @{ D = "34.05";      E = "  -0.60%  " }
    30 = @{ D = "9.20";       E = "  +0.08%  " }
    31 = @{ E = "  -6.02%  " }
    32 = @{ D = "5.00";       E = "  +7.99%  " }
    33 = @{ E = "  +2.56%  " }
    34 = @{ D = "0.999";      E = "  -0.03%  " }
    35 = @{ D = "18.09";      E = "  +5.50%  " }
    36 = @{ E = "  -1.51%  " }
    37 = @{ E = "  +0.74%  " }
    38 = @{ E = "  +1.36%  " }
    39 = @{ E = "  -0.14%  " }
    40 = @{ D = "2.79";       E = "  -1.26%  " }
    41 = @{ E = "  +0.01%  " }
    42 = @{ D = "1.994.23";   E = "  -0.34%  " }
    43 = @{ E = "  +1.11%  " }
    44 = @{ E = "  -5.93%  " }
    45 = @{ D = "10.24";      E = "  +0.97%  " }
    46 = @{ D = "17.69";      E = "  -0.56%  " }
    47 = @{ E = "  -0.29%  " }
    48 = @{ B = "MultiversX"; C = "https://coinranking.com/coin/omwkOTglq+multiversx-egld";      D = "55.11"; E = "  -1.69%  " }
    49 = @{ B = "BitcoinSV";  C = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv";     D = "76.13"; E = "  +8.15%  " }
    50 = @{ D = "2.546.73";   E = "  +0.45%  " }
    51 = @{ E = "  +0.47%  " }
}

foreach ($row in $rowUpdates.Keys) {
    $fields = $rowUpdates[$row]

    if ($fields.Contains("B")) {
        $ws.Range("B$row").Value = $fields["B"]
    }
    if ($fields.Contains("C")) {
        $ws.Range("C$row").Value = $fields["C"]
    }
    if ($fields.Contains("D")) {
        Set-TextValue $ws.Range("D$row") $fields["D"]
    }
    if ($fields.Contains("E")) {
        Set-TextValue $ws.Range("E$row") $fields["E"]
    }
}
